$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks price history: each column from B onward is a
# timestamped snapshot, and the last two columns hold the product name
# ("nom") and its URL ("url_produit"). A new price check was performed,
# so a new snapshot column is inserted right before the "nom" column
# (currently AM), pushing "nom" to AN and "url_produit" to AO.
$ws.Columns("AM").Insert()

# Header for the newly inserted snapshot column: the timestamp of this run.
$ws.Range("AM1").Value = "2026-01-29 10:25:55"

# Figure out how many data rows exist (row 1 is the header).
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 206 }

# Populate the new snapshot column: carry forward the most recent known
# price (previously in column AL) as today's price for rows that had
# one; rows with no recorded price yet stay blank, same as before.
for ($r = 2; $r -le $lastRow; $r++) {
    $alCell = $ws.Range("AL" + $r)
    $alVal = $alCell.Value2
    if ($alVal -ne $null -and $alVal -ne "") {
        $ws.Range("AM" + $r).Value = $alVal
    }
}
